$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.713.74'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.392.87'
$ws.Range("E3").Value = '  +2.48%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.30%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").Value = '2.393.62'
$ws.Range("E9").Value = '  +2.51%  '
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("E15").Value = '  +9.09%  '
$ws.Range("D16").Value = '2.825.28'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("D17").Value = '61.591.72'
$ws.Range("E17").Value = '  +1.86%  '
$ws.Range("D18").Value = '2.388.56'
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("E19").Value = '  +5.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '322.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.79%  '
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -5.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '555.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '2.500.68'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.44%  '
$ws.Range("D31").Value = '0.0₃0916'
$ws.Range("E31").Value = '  +3.19%  '
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("E36").Value = '  +10.46%  '
$ws.Range("E38").Value = '  +8.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.381'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '145.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.12%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.91%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.38%  '
$ws.Range("E46").Value = '  +2.92%  '
$ws.Range("E47").Value = '  +4.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.586'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0904'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("E51").Value = '  +1.64%  '
